$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1937984496124031
$ws.Range("C2").Value = 0.5852713178294574
$ws.Range("J2").Value = 0.01162790697674419
$ws.Range("P2").Value = 0.1317829457364341
$ws.Range("S2").Value = 0.07751937984496124
$ws.Range("J3").Value = 0.01333333333333333
$ws.Range("P3").Value = 0.8133333333333334
$ws.Range("S3").Value = 0.1733333333333333
$ws.Range("J4").Value = 0.06060606060606061
$ws.Range("P4").Value = 0.5757575757575758
$ws.Range("S4").Value = 0.3636363636363636
$ws.Range("B6").Value = 0.04
$ws.Range("D6").Value = 0.02
$ws.Range("F6").Value = 0.1
$ws.Range("J6").Value = 0.228
$ws.Range("O6").Value = 0.02
$ws.Range("Q6").Value = 0.168
$ws.Range("R6").Value = 0.092
$ws.Range("S6").Value = 0.332
$ws.Range("B7").Value = 0.09302325581395349
$ws.Range("E7").Value = 0.004651162790697674
$ws.Range("F7").Value = 0.07906976744186046
$ws.Range("J7").Value = 0.1441860465116279
$ws.Range("O7").Value = 0.04651162790697674
$ws.Range("Q7").Value = 0.1302325581395349
$ws.Range("R7").Value = 0.1302325581395349
$ws.Range("S7").Value = 0.3720930232558139
$ws.Range("B8").Value = 0.07633587786259542
$ws.Range("D8").Value = 0.02035623409669211
$ws.Range("E8").Value = 0.002544529262086514
$ws.Range("F8").Value = 0.06361323155216285
$ws.Range("J8").Value = 0.1399491094147583
$ws.Range("O8").Value = 0.02290076335877863
$ws.Range("Q8").Value = 0.1577608142493639
$ws.Range("R8").Value = 0.1424936386768448
$ws.Range("S8").Value = 0.3740458015267176
$ws.Range("B9").Value = 0.09210526315789473
$ws.Range("D9").Value = 0.008771929824561403
$ws.Range("F9").Value = 0.08771929824561403
$ws.Range("J9").Value = 0.1403508771929824
$ws.Range("O9").Value = 0.02192982456140351
$ws.Range("Q9").Value = 0.1491228070175439
$ws.Range("R9").Value = 0.1403508771929824
$ws.Range("S9").Value = 0.3596491228070176
$ws.Range("B10").Value = 0.08996282527881042
$ws.Range("D10").Value = 0.01412639405204461
$ws.Range("E10").Value = 0.0007434944237918215
$ws.Range("F10").Value = 0.07732342007434945
$ws.Range("J10").Value = 0.104089219330855
$ws.Range("O10").Value = 0.02156133828996282
$ws.Range("Q10").Value = 0.2133828996282528
$ws.Range("R10").Value = 0.1152416356877323
$ws.Range("S10").Value = 0.3635687732342007
$ws.Range("G11").Value = 0.1337386018237082
$ws.Range("J11").Value = 0.06990881458966565
$ws.Range("K11").Value = 0.1914893617021277
$ws.Range("L11").Value = 0.5927051671732523
$ws.Range("S11").Value = 0.0121580547112462
$ws.Range("G12").Value = 0.7355769230769231
$ws.Range("J12").Value = 0.1875
$ws.Range("K12").Value = 0.009615384615384616
$ws.Range("L12").Value = 0.04326923076923077
$ws.Range("S12").Value = 0.02403846153846154
$ws.Range("G13").Value = 0.6944444444444444
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.05555555555555555
$ws.Range("F15").Value = 0.04489795918367347
$ws.Range("H15").Value = 0.1551020408163265
$ws.Range("I15").Value = 0.06122448979591837
$ws.Range("J15").Value = 0.2897959183673469
$ws.Range("K15").Value = 0.06938775510204082
$ws.Range("M15").Value = 0.004081632653061225
$ws.Range("O15").Value = 0.07346938775510205
$ws.Range("S15").Value = 0.3020408163265306
$ws.Range("F16").Value = 0.01169590643274854
$ws.Range("H16").Value = 0.1637426900584795
$ws.Range("I16").Value = 0.05847953216374269
$ws.Range("J16").Value = 0.4327485380116959
$ws.Range("K16").Value = 0.1052631578947368
$ws.Range("M16").Value = 0.005847953216374269
$ws.Range("O16").Value = 0.03508771929824561
$ws.Range("S16").Value = 0.1871345029239766
$ws.Range("F17").Value = 0.0111358574610245
$ws.Range("H17").Value = 0.2138084632516704
$ws.Range("I17").Value = 0.0957683741648107
$ws.Range("J17").Value = 0.4142538975501114
$ws.Range("K17").Value = 0.09799554565701558
$ws.Range("M17").Value = 0.0111358574610245
$ws.Range("N17").Value = 0.0022271714922049
$ws.Range("O17").Value = 0.0400890868596882
$ws.Range("S17").Value = 0.1135857461024499
$ws.Range("F18").Value = 0.01030927835051546
$ws.Range("H18").Value = 0.154639175257732
$ws.Range("I18").Value = 0.1099656357388316
$ws.Range("J18").Value = 0.4398625429553265
$ws.Range("K18").Value = 0.1030927835051546
$ws.Range("M18").Value = 0.02061855670103093
$ws.Range("O18").Value = 0.05154639175257732
$ws.Range("S18").Value = 0.1099656357388316
$ws.Range("F19").Value = 0.01364365971107544
$ws.Range("H19").Value = 0.1532905296950241
$ws.Range("I19").Value = 0.1043338683788122
$ws.Range("J19").Value = 0.4044943820224719
$ws.Range("K19").Value = 0.122792937399679
$ws.Range("M19").Value = 0.02006420545746388
$ws.Range("O19").Value = 0.07945425361155699
$ws.Range("S19").Value = 0.1019261637239165
